# "RF classify (read data header bug)"
#
# The sheet used to contain a duplicate "raw header" table in columns H:M
# (start/end time read back from a video-timestamp log) that A:F's
# start-time/end-time columns referenced via formulas (e.g. A2 = H2-B$9,
# B2 = I2-B$9). That raw/duplicate block, together with the "Video start
# time:" helper row (row 9) it depended on, was buggy scaffolding and is
# removed here. Because the formulas in A:B depended on that data, they
# are first frozen to their already-computed static values so the
# numbers in A:F are preserved once the source columns/row disappear.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Freeze the formulas in A2:B7 (they reference columns H/I and B$9)
#    to plain values, since their source data is about to be deleted.
$timeRange = $ws.Range("A2:B7")
$timeRange.Value = $timeRange.Value2

# 2) Remove the duplicate/raw header-and-data block in columns H:M.
$ws.Range("H1:M7").EntireColumn.Delete()

# 3) Remove row 9 ("Video start time:" label and its value). This shifts
#    the lone remaining helper row (old row 11, the 1-second nudge used
#    by the old H/I formulas) up to row 10.
$ws.Range("A9").EntireRow.Delete()

# 4) Leave the selection where it would land after clicking the header
#    of the now-empty column immediately to the right of the data.
$ws.Range("H1:Q1048576").Select()
